# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells in AD1:AF1 - copy the existing header formatting (style)
# from AC1, then set the header text.
$headerFormat = $ws.Range("AC1")
$headerFormat.Copy()

$wins = $ws.Range("AD1")
$wins.PasteSpecial(-4122)
$wins.Value = "Wins"

$losses = $ws.Range("AE1")
$losses.PasteSpecial(-4122)
$losses.Value = "Losses"

$ties = $ws.Range("AF1")
$ties.PasteSpecial(-4122)
$ties.Value = "Ties"

# Every team row (2-58) shares the same record: 69 wins, 93 losses, 0 ties.
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 69
    $ws.Cells.Item($r, 31).Value = 93
    $ws.Cells.Item($r, 32).Value = 0
}
